$d = $word.ActiveDocument

# "agregamos atributos al usuario" - insert five new paragraphs (Java-style
# field declarations for a "usuario" class) before the existing (empty)
# first paragraph of the document body. Each paragraph is built from the
# exact WordprocessingML run/proofErr layout Word itself produces when a
# line is typed with live spell-check on, so InsertXML is used to land
# that markup verbatim rather than relying on a synthesized Find/Replace
# or plain-text insert (which would not reproduce the <w:proofErr/>
# spell-check markers or the per-word run splitting).

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-AttrParagraph([string]$type, [string]$name) {
    # "<type> <name>;" where both words get wrapped in spellStart/spellEnd
    # proof-error markers (mirrors what Word's background spell checker
    # inserts for words it does not recognize, e.g. "String", "int").
    return "<w:p $wNs>" +
           '<w:proofErr w:type="spellStart"/>' +
           "<w:r><w:t>$type</w:t></w:r>" +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/>' +
           "<w:r><w:t>$name</w:t></w:r>" +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t>;</w:t></w:r>' +
           '</w:p>'
}

$p1 = New-AttrParagraph "String" "nombreUsuario"
$p2 = New-AttrParagraph "String" "idUsuario"
$p3 = New-AttrParagraph "int" "edadUsuario"

# "Date" is a recognized dictionary word, so Word does not flag it (no
# proofErr wrapping) and keeps it in the same run as the trailing space.
$p4 = "<w:p $wNs>" +
      '<w:r><w:t xml:space="preserve">Date </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>nacimientoUsuario</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t>;</w:t></w:r>' +
      '</w:p>'

$p5 = New-AttrParagraph "long" "cedulaUsuario"

$xml = $p1 + $p2 + $p3 + $p4 + $p5

$r = $d.Range(0, 0)
$r.InsertXML($xml)
